# Auto-generated Excel COM-interop script
# Applies: (1) timestamp update, (2) country name reshuffle (58 rows),
# (3) 113 numeric cell updates across the COVID data table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Update the "datos actualizados" timestamp in A1 ---
$ws.Range("A1").Value = "Datos actualizados a 23 de Marzo de 2020 a las 18:16"

# --- 2) Re-shuffle country names in column A (rows 4-197) ---
# Values are assigned per-row to their new target country text;
# the engine dedupes/reorders the shared-string table on save.
$countryUpdates = @{
    26 = "Chequia"
    27 = "Turquia"
    58 = "Colombia"
    59 = "Libano"
    60 = "Argentina"
    61 = "Irak"
    62 = "Serbia"
    63 = "Republica Dominicana"
    107 = "Liechtenstein"
    108 = "Uzbekistan"
    110 = "Cuba"
    111 = "Afganistan"
    112 = "Consejo Danes para los Refugiados"
    113 = "Nigeria"
    121 = "Macao"
    122 = "Ghana"
    129 = "Guatemala"
    130 = "Guyana"
    134 = "Islas Virgenes de los Estados Unidos"
    135 = "Barbados"
    150 = "Haiti"
    151 = "San Martin (Parte Francesa)"
    155 = "Suazilandia"
    157 = "Guinea"
    159 = "Cabo Verde"
    160 = "Fiyi"
    161 = "Namibia"
    162 = "San Bartolome"
    163 = "Zambia"
    164 = "Liberia"
    165 = "Congo"
    166 = "El Salvador"
    167 = "Republica de Africa Central"
    170 = "Mauritania"
    171 = "San Martin (Parte Holandesa)"
    172 = "Nicaragua"
    173 = "Angola"
    174 = "Niger"
    175 = "Benin"
    176 = "Butan"
    177 = "Santa Lucia"
    180 = "Sudan"
    181 = "Belice"
    182 = "San Vicente y las Granadinas"
    183 = "Uganda"
    184 = "Republica de Yibuti"
    185 = "Dominica"
    186 = "Republica del Chad"
    187 = "Eritrea"
    189 = "Timor Oriental"
    190 = "Somalia"
    191 = "Santa Sede"
    192 = "Papua Nueva Guinea"
    193 = "Antigua y Barbuda"
    194 = "Mozambique"
    195 = "Siria"
    196 = "Islas Turcas y Caicos"
    197 = "Granada"
}
foreach ($r in $countryUpdates.Keys) {
    $ws.Cells.Item([int]$r, 1).Value = $countryUpdates[$r]
}

# --- 3) Apply numeric value updates to columns B-H ---
$colIndex = @{ "B" = 2; "C" = 3; "D" = 4; "E" = 5; "F" = 6; "G" = 7; "H" = 8 }
$cellUpdates = @(
    @(5, "F", 3009),
    @(6, "B", 40841),
    @(6, "C", 7295),
    @(6, "E", 40171),
    @(6, "G", 64),
    @(6, "H", 483),
    @(8, "B", 28798),
    @(8, "C", 3925),
    @(8, "E", 28260),
    @(8, "G", 22),
    @(8, "H", 116),
    @(21, "B", 1629),
    @(21, "C", 83),
    @(21, "E", 1602),
    @(26, "C", 116),
    @(26, "D", 6),
    @(26, "E", 1229),
    @(26, "F", 19),
    @(26, "H", 1),
    @(27, "B", 1236),
    @(27, "C", 0),
    @(27, "D", 0),
    @(27, "E", 1206),
    @(27, "F", 0),
    @(27, "H", 30),
    @(58, "B", 277),
    @(58, "C", 46),
    @(58, "D", 3),
    @(58, "E", 271),
    @(58, "F", 0),
    @(58, "G", 1),
    @(58, "H", 3),
    @(59, "B", 267),
    @(59, "C", 19),
    @(59, "D", 8),
    @(59, "E", 255),
    @(59, "F", 4),
    @(60, "C", 0),
    @(60, "D", 27),
    @(60, "E", 235),
    @(60, "G", 0),
    @(60, "H", 4),
    @(61, "B", 266),
    @(61, "C", 33),
    @(61, "D", 62),
    @(61, "E", 181),
    @(61, "F", 0),
    @(61, "G", 3),
    @(61, "H", 23),
    @(62, "B", 249),
    @(62, "C", 27),
    @(62, "D", 3),
    @(62, "E", 244),
    @(62, "F", 12),
    @(62, "H", 2),
    @(63, "B", 245),
    @(63, "C", 43),
    @(63, "D", 0),
    @(63, "E", 242),
    @(63, "G", 0),
    @(90, "B", 95),
    @(90, "C", 13),
    @(90, "D", 2),
    @(90, "E", 93),
    @(93, "B", 87),
    @(93, "C", 3),
    @(93, "E", 85),
    @(107, "C", 9),
    @(108, "C", 3),
    @(110, "C", 5),
    @(110, "D", 0),
    @(110, "E", 39),
    @(110, "F", 3),
    @(111, "B", 40),
    @(111, "C", 0),
    @(111, "D", 1),
    @(111, "E", 38),
    @(112, "D", 0),
    @(112, "E", 35),
    @(112, "G", 0),
    @(113, "B", 36),
    @(113, "C", 6),
    @(113, "D", 2),
    @(113, "E", 33),
    @(113, "G", 1),
    @(121, "B", 25),
    @(121, "C", 3),
    @(121, "D", 10),
    @(121, "E", 15),
    @(121, "H", 0),
    @(122, "C", 1),
    @(122, "D", 0),
    @(122, "E", 23),
    @(122, "H", 1),
    @(134, "C", 11),
    @(135, "C", 3),
    @(150, "B", 6),
    @(150, "C", 4),
    @(150, "E", 6),
    @(151, "C", 0),
    @(155, "C", 0),
    @(157, "C", 2),
    @(160, "C", 1),
    @(163, "C", 0),
    @(171, "C", 1),
    @(177, "E", 2),
    @(177, "H", 0),
    @(180, "B", 2),
    @(180, "C", 0),
    @(180, "H", 1),
    @(181, "C", 1),
    @(194, "C", 0),
    @(196, "C", 1)
)
foreach ($upd in $cellUpdates) {
    $r = $upd[0]
    $c = $colIndex[$upd[1]]
    $v = $upd[2]
    $ws.Cells.Item($r, $c).Value = $v
}

Write-Host "Applied timestamp + $($countryUpdates.Count) country reorders + $($cellUpdates.Count) numeric updates"
